$wb = $excel.ActiveWorkbook

# Rename sheets (new timestamp-based task-order identifiers)
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168753160181"
$wb.Worksheets.Item(2).Name = "NB_TO-16511687563795125"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687563805163"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687564422789"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511687565190196"

# Sheet 1 (GNG_TO): update stim filenames
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511687531313653.csv"
$ws1.Range("B3").Value = "GNG_stims-16511687531441283.csv"
$ws1.Range("B4").Value = "go_stims-165116875314513.csv"
$ws1.Range("B5").Value = "GNG_stims-1651168753159213.csv"

# Sheet 2 (NB_TO): update stim filenames
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511687548141022.csv"
$ws2.Range("B3").Value = "ZB-match_2-16511687535823379.csv"
$ws2.Range("B4").Value = "TB-16511687551242456.csv"
$ws2.Range("B5").Value = "ZB-match_3-1651168753177839.csv"
$ws2.Range("B6").Value = "TB-16511687563565955.csv"
$ws2.Range("B7").Value = "TB-16511687554338725.csv"
$ws2.Range("B8").Value = "ZB-match_5-16511687533375025.csv"
$ws2.Range("B9").Value = "OB-1651168753792312.csv"
$ws2.Range("B10").Value = "OB-16511687541482425.csv"

# Sheet 3 (RS_TO): swap eyes open / eyes closed
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL_TO): update stim filenames
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-1651168756394822.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168756382517.csv"
$ws4.Range("B4").Value = "MM_stims-16511687564259543.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687563958127.csv"
$ws4.Range("B6").Value = "MM_stims-1651168756441289.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687564269211.csv"

# Sheet 5 (vSAT_TO): update stim filenames
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-1651168756449077.csv"
$ws5.Range("B3").Value = "SAT_stims-16511687564737713.csv"
$ws5.Range("B4").Value = "vSAT_stims-16511687565038464.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511687564880314.csv"
